$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure numeric-looking text values (prices/percentages) are stored as text,
# matching the source data which uses "." as a thousands separator in some rows.
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "21.090.03"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  -4.29%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.506.97"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  -2.99%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.006"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  +0.50%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "1.006"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  +0.50%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "284.32"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  -2.02%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.3846"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  -2.42%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3126"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  -2.97%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "42.59"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  -2.88%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.06973"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  -3.68%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "1.039"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  -3.30%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.006"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  +0.52%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "5.598"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  -1.15%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "17.81"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  -4.99%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "1.512.51"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  -2.51%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "6.365"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  -3.84%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.00001068"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  -5.58%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.06567"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  -0.05%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "81.79"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  -1.91%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "1.006"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  +0.54%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "5.979"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  -4.59%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "15.11"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  -2.52%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "10.83"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  -4.09%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.346"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  -0.93%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "21.086.93"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  -4.36%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.347"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  -2.53%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "147.75"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  -0.84%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "18.00"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  -2.81%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "4.789"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  -1.83%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.681.01"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  -2.61%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "114.45"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  -3.51%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "5.890"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  +1.17%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.9520"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  -2.20%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.07935"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  -4.85%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "8.390"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  -7.48%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "5.060"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  -0.89%  "

$ws.Range("B37").Value = "WEMIXTOKEN"
$ws.Range("C37").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.469"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  -8.28%  "

$ws.Range("B38").Value = "Aptos"
$ws.Range("C38").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "11.24"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  +5.50%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.05781"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  -3.75%  "

$ws.Range("B40").Value = "Frax"
$ws.Range("C40").Value = "https://coinranking.com/coin/KfWtaeV1W+frax-frax"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.006"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  +0.53%  "

$ws.Range("B41").Value = "VeChain"
$ws.Range("C41").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.02120"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  -6.13%  "

$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  -4.18%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.1971"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  -3.06%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.5613"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  -3.21%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "12.93"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  -0.86%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "3.668"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  -2.01%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.5423"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  -2.53%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.137"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  +0.42%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.844"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  -2.68%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "113.62"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  -3.44%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.06550"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  -3.96%  "
